$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 57570.43
$ws.Range("I21").Value = 80019
$ws.Range("J21").Value = 40734
$ws.Range("K21").Value = 80019
$ws.Range("L21").Value = 40734
$ws.Range("M21").Value = -79551
$ws.Range("N21").Value = -41670
$ws.Range("H23").Value = 57570.43
$ws.Range("I23").Value = 80019
$ws.Range("J23").Value = 40734
$ws.Range("K23").Value = 80019
$ws.Range("L23").Value = 40734
$ws.Range("M23").Value = -79785
$ws.Range("N23").Value = -41202
$ws.Range("H29").Value = 2100.1667
$ws.Range("I29").Value = 899.25
$ws.Range("J29").Value = 4502
$ws.Range("K29").Value = 2697.75
$ws.Range("L29").Value = 13506
$ws.Range("M29").Value = -2416.75
$ws.Range("N29").Value = -14068
$ws.Range("H38").Value = 3690.6924
$ws.Range("I38").Value = 109.875
$ws.Range("J38").Value = 9420
$ws.Range("K38").Value = 329.625
$ws.Range("L38").Value = 28260
$ws.Range("M38").Value = 42.375
$ws.Range("N38").Value = -29004
$ws.Range("H44").Value = 18571.428
$ws.Range("J44").Value = 18571.428
$ws.Range("L44").Value = 18571.428
$ws.Range("N44").Value = -19495.428
$ws.Range("H58").Value = 7322.4287
$ws.Range("I58").Value = 946
$ws.Range("J58").Value = 18800
$ws.Range("K58").Value = 2838
$ws.Range("L58").Value = 56400
$ws.Range("M58").Value = -2688
$ws.Range("N58").Value = -56700
$ws.Range("H87").Value = 23859.54
$ws.Range("J87").Value = 23859.54
$ws.Range("L87").Value = 23859.54
$ws.Range("N87").Value = -26355.54
$ws.Range("H90").Value = 23859.54
$ws.Range("J90").Value = 23859.54
$ws.Range("L90").Value = 71578.62
$ws.Range("N90").Value = -84058.62
$ws.Range("H112").Value = 1564.8877
$ws.Range("I112").Value = 700
$ws.Range("J112").Value = 1573.8041
$ws.Range("K112").Value = 2100
$ws.Range("L112").Value = 4721.4123
$ws.Range("M112").Value = -992
$ws.Range("N112").Value = -6937.4123
$ws.Range("H116").Value = 10323.154
$ws.Range("I116").Value = 2065.8333
$ws.Range("J116").Value = 17400.857
$ws.Range("K116").Value = 2065.8333
$ws.Range("L116").Value = 17400.857
$ws.Range("M116").Value = 1376.1667
$ws.Range("N116").Value = -24284.857
$ws.Range("H129").Value = 1114.6428
$ws.Range("I129").Value = 350
$ws.Range("J129").Value = 1133.2927
$ws.Range("K129").Value = 1050
$ws.Range("L129").Value = 3399.8781
$ws.Range("M129").Value = 3950
$ws.Range("N129").Value = -13399.8781
$ws.Range("H137").Value = 2617.5945
$ws.Range("I137").Value = 1879.7037
$ws.Range("K137").Value = 5639.1111
$ws.Range("M137").Value = -3089.1111
$ws.Range("H138").Value = 2474.52
$ws.Range("I138").Value = 1107.5
$ws.Range("J138").Value = 2860.0898
$ws.Range("K138").Value = 3322.5
$ws.Range("L138").Value = 8580.269400000001
$ws.Range("M138").Value = 1817.5
$ws.Range("N138").Value = -18860.2694

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 5008723
$ws.Range("I6").Value = 6674310
$ws.Range("J6").Value = 11961.5
$ws.Range("K6").Value = 6674310
$ws.Range("L6").Value = 11961.5
$ws.Range("M6").Value = -6674137
$ws.Range("N6").Value = -12307.5
$ws.Range("H32").Value = 7552.254
$ws.Range("I32").Value = 4745.2
$ws.Range("K32").Value = 4745.2
$ws.Range("M32").Value = -4458.2
$ws.Range("H63").Value = 5774204
$ws.Range("I63").Value = 10656600
$ws.Range("J63").Value = 4100
$ws.Range("K63").Value = 10656600
$ws.Range("L63").Value = 4100
$ws.Range("M63").Value = -10655914
$ws.Range("N63").Value = -5472
$ws.Range("H66").Value = 5774204
$ws.Range("I66").Value = 10656600
$ws.Range("J66").Value = 4100
$ws.Range("K66").Value = 53283000
$ws.Range("L66").Value = 20500
$ws.Range("M66").Value = -53279568
$ws.Range("N66").Value = -27364
$ws.Range("H132").Value = 3214.9412
$ws.Range("I132").Value = 2400.9048
$ws.Range("J132").Value = 4529.923
$ws.Range("K132").Value = 7202.714399999999
$ws.Range("L132").Value = 13589.769
$ws.Range("M132").Value = -4672.714399999999
$ws.Range("N132").Value = -18649.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 67466.664
$ws.Range("J70").Value = 67466.664
$ws.Range("L70").Value = 67466.664
$ws.Range("N70").Value = -68052.664
$ws.Range("H73").Value = 67466.664
$ws.Range("J73").Value = 67466.664
$ws.Range("L73").Value = 67466.664
$ws.Range("N73").Value = -69494.664
$ws.Range("H82").Value = 22582.7
$ws.Range("I82").Value = 6692.75
$ws.Range("J82").Value = 33176
$ws.Range("K82").Value = 6692.75
$ws.Range("L82").Value = 33176
$ws.Range("M82").Value = -6309.75
$ws.Range("N82").Value = -33942
$ws.Range("H85").Value = 22582.7
$ws.Range("I85").Value = 6692.75
$ws.Range("J85").Value = 33176
$ws.Range("K85").Value = 6692.75
$ws.Range("L85").Value = 33176
$ws.Range("M85").Value = -5366.75
$ws.Range("N85").Value = -35828

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35719676
$ws.Range("I31").Value = 2900
$ws.Range("K31").Value = 2900
$ws.Range("M31").Value = -2605
$ws.Range("H34").Value = 35719676
$ws.Range("I34").Value = 2900
$ws.Range("K34").Value = 2900
$ws.Range("M34").Value = -2698
$ws.Range("H41").Value = 36824.8
$ws.Range("J41").Value = 45016.25
$ws.Range("L41").Value = 45016.25
$ws.Range("N41").Value = -45872.25
$ws.Range("H50").Value = 26893.125
$ws.Range("J50").Value = 26893.125
$ws.Range("L50").Value = 26893.125
$ws.Range("N50").Value = -28143.125
$ws.Range("H59").Value = 29178.555
$ws.Range("J59").Value = 29178.555
$ws.Range("L59").Value = 29178.555
$ws.Range("N59").Value = -31468.555
$ws.Range("H60").Value = 29404.385
$ws.Range("I60").Value = 3000
$ws.Range("J60").Value = 31604.75
$ws.Range("K60").Value = 3000
$ws.Range("L60").Value = 31604.75
$ws.Range("M60").Value = -2489
$ws.Range("N60").Value = -32626.75
$ws.Range("H74").Value = 49828.5
$ws.Range("J74").Value = 49828.5
$ws.Range("L74").Value = 49828.5
$ws.Range("N74").Value = -51576.5
$ws.Range("H77").Value = 49828.5
$ws.Range("J77").Value = 49828.5
$ws.Range("L77").Value = 149485.5
$ws.Range("N77").Value = -158221.5
$ws.Range("H87").Value = 21077.777
$ws.Range("J87").Value = 21077.777
$ws.Range("L87").Value = 21077.777
$ws.Range("N87").Value = -23449.777
$ws.Range("H90").Value = 21077.777
$ws.Range("J90").Value = 21077.777
$ws.Range("L90").Value = 63233.33099999999
$ws.Range("N90").Value = -75089.33099999999
$ws.Range("H132").Value = 3142.1765
$ws.Range("I132").Value = 1478.9412
$ws.Range("J132").Value = 4805.4116
$ws.Range("K132").Value = 4436.8236
$ws.Range("L132").Value = 14416.2348
$ws.Range("M132").Value = -1906.8236
$ws.Range("N132").Value = -19476.2348

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 9914.143
$ws.Range("I26").Value = 30124.75
$ws.Range("J26").Value = 1829.9
$ws.Range("K26").Value = 90374.25
$ws.Range("L26").Value = 5489.700000000001
$ws.Range("M26").Value = -90086.25
$ws.Range("N26").Value = -6065.700000000001
$ws.Range("H113").Value = 567.64703
$ws.Range("I113").Value = 564.8276
$ws.Range("J113").Value = 571.36365
$ws.Range("K113").Value = 1694.4828
$ws.Range("L113").Value = 1714.09095
$ws.Range("M113").Value = 475.5172000000002
$ws.Range("N113").Value = -6054.09095
$ws.Range("H121").Value = 2533.712
$ws.Range("I121").Value = 256.5
$ws.Range("J121").Value = 2791.5095
$ws.Range("K121").Value = 769.5
$ws.Range("L121").Value = 8374.5285
$ws.Range("M121").Value = 540.5
$ws.Range("N121").Value = -10994.5285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 15627545
$ws.Range("I80").Value = 22729720
$ws.Range("J80").Value = 2761.2
$ws.Range("K80").Value = 22729720
$ws.Range("L80").Value = 2761.2
$ws.Range("M80").Value = -22728722
$ws.Range("N80").Value = -4757.2
$ws.Range("H83").Value = 15627545
$ws.Range("I83").Value = 22729720
$ws.Range("J83").Value = 2761.2
$ws.Range("K83").Value = 113648600
$ws.Range("L83").Value = 13806
$ws.Range("M83").Value = -113643608
$ws.Range("N83").Value = -23790
$ws.Range("H132").Value = 3014.1155
$ws.Range("I132").Value = 1804.7059
$ws.Range("J132").Value = 5298.5557
$ws.Range("K132").Value = 5414.1177
$ws.Range("L132").Value = 15895.6671
$ws.Range("M132").Value = -2884.1177
$ws.Range("N132").Value = -20955.6671
$ws.Range("H133").Value = 44266.9
$ws.Range("J133").Value = 44266.9
$ws.Range("L133").Value = 44266.9
$ws.Range("N133").Value = -54386.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14495176
$ws.Range("I132").Value = 909.36365
$ws.Range("J132").Value = 27781588
$ws.Range("K132").Value = 2728.09095
$ws.Range("L132").Value = 83344764
$ws.Range("M132").Value = -198.0909499999998
$ws.Range("N132").Value = -83349824
